$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price/volume updates ---
# D-column values are forced to Text before assignment (and the
# temporary Text number-format is cleared afterwards) so that
# numeric-looking strings like "0.999" stay literal text cells
# instead of being auto-coerced into numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '71.180.12'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.29%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.583.92'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -4.87%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.06'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -3.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.32'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.11%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.513'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.585.13'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -4.71%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.166'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.59%  '
$ws.Range("E11").Value = '  +0.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.352'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.76%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.86'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.42%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.063.27'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -4.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000184'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '71.046.46'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.30'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -4.06%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.595.90'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -4.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.86'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.85%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.70'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -4.97%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '365.69'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.54%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.00'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -4.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.99'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.15'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -4.62%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.29'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -6.25%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.763.25'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.12%  '
$ws.Range("E29").Value = '  +0.45%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0930'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -5.79%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.85'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -3.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.78'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.44%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '158.45'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -3.84%  '
$ws.Range("E37").Value = '  +5.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.95'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.99%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.86'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.34'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -3.84%  '
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("E42").Value = '  -6.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.51'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.54%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.78'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -5.72%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.321'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -4.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.64'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.86%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '147.34'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -6.57%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.59'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -4.51%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.533'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -5.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.65'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -7.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.597'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.15%  '

# --- Row 24/25 swap: Litecoin <-> Dai with updated price/volume ---
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.96'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.66%  '

# --- Row 32/33 swap: Bittensor <-> Fetch.AI with updated price/volume ---
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.31'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.35%  '
$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '481.64'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -5.84%  '
